$d = $word.ActiveDocument

# Directeur technique : 8160€ -> 4760€
$d.Content.Find.Execute("8160", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4760", 2)

# Chef de projet : 6100€ -> 12200€
$d.Content.Find.Execute("6100", $true, $false, $false, $false, $false,
                         $true, 1, $false, "12200", 2)

# Montant total estimé : 14260€ -> 16960€
$d.Content.Find.Execute("14260", $true, $false, $false, $false, $false,
                         $true, 1, $false, "16960", 2)

# Montant total : 56805€ -> 59505€
$d.Content.Find.Execute("56805", $true, $false, $false, $false, $false,
                         $true, 1, $false, "59505", 2)
